$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 - this shifts existing rows 23-61 down to 24-62,
# matching the rest of the weekly data already present.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with this week's record.
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44533
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112021
$ws.Range("G23").Value = "Ají"
$ws.Range("H23").Value = "Inferno"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = 18000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 18800
$ws.Range("N23").Value = "$/caja 15 kilos"
$ws.Range("O23").Value = "Provincia de Huasco"
$ws.Range("P23").Value = 1253
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = "Hortaliza"

# D column (Fecha) needs the date number format used throughout the sheet.
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
